$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.640.64'
$ws.Range('E2').Value = '  +3.43%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.492.13'
$ws.Range('E3').Value = '  +2.13%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.93'
$ws.Range('E5').Value = '  +2.17%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.30'
$ws.Range('E6').Value = '  +3.31%  '

# Row 7
$ws.Range('E7').Value = '  -0.06%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.609'
$ws.Range('E8').Value = '  +12.07%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.499.67'
$ws.Range('E9').Value = '  +2.33%  '

# Row 10
$ws.Range('E10').Value = '  -1.46%  '

# Row 11
$ws.Range('E11').Value = '  +2.76%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.449'
$ws.Range('E12').Value = '  +2.93%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.098.18'
$ws.Range('E13').Value = '  +2.29%  '

# Row 14
$ws.Range('E14').Value = '  +0.26%  '

# Row 15
$ws.Range('E15').Value = '  +2.47%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '28.80'
$ws.Range('E16').Value = '  +6.35%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.643.68'
$ws.Range('E17').Value = '  +3.15%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.476.66'
$ws.Range('E18').Value = '  +2.16%  '

# Row 19
$ws.Range('E19').Value = '  +2.97%  '

# Row 20
$ws.Range('E20').Value = '  +1.37%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '389.03'
$ws.Range('E21').Value = '  +0.64%  '

# Row 22
$ws.Range('E22').Value = '  +0.86%  '

# Row 23
$ws.Range('E23').Value = '  +3.56%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.43'
$ws.Range('E24').Value = '  +1.84%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  +0.25%  '

# Row 26
$ws.Range('E26').Value = '  +5.14%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.28'
$ws.Range('E27').Value = '  +8.25%  '

# Row 28
$ws.Range('E28').Value = '  +1.62%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.48%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.30'
$ws.Range('E30').Value = '  +4.86%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.46'
$ws.Range('E31').Value = '  +7.01%  '

# Row 32
$ws.Range('E32').Value = '  +3.22%  '

# Row 33
$ws.Range('E33').Value = '  +2.00%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.27'
$ws.Range('E34').Value = '  +6.23%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.59'
$ws.Range('E35').Value = '  +8.35%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.83'
$ws.Range('E36').Value = '  +2.29%  '

# Row 37
$ws.Range('E37').Value = '  +6.64%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.046.04'
$ws.Range('E38').Value = '  +5.24%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0778'
$ws.Range('E39').Value = '  +1.63%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '27.37'
$ws.Range('E40').Value = '  +2.00%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0326'
$ws.Range('E41').Value = '  +2.43%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.59'
$ws.Range('E42').Value = '  +4.63%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '43.14'
$ws.Range('E43').Value = '  +5.27%  '

# Row 44
$ws.Range('E44').Value = '  +1.61%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.781'
$ws.Range('E45').Value = '  +2.14%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.97'
$ws.Range('E46').Value = '  +10.07%  '

# Row 47
$ws.Range('E47').Value = '  +3.82%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '319.72'
$ws.Range('E48').Value = '  +9.91%  '

# Row 49
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.78'
$ws.Range('E49').Value = '  +5.17%  '

# Row 50
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.23'
$ws.Range('E50').Value = '  +1.64%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.110'
$ws.Range('E51').Value = '  +6.36%  '
